$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Mon Feb 24 23:06:58 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:07:13 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:07:27 EST 2025"
